$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 31 ---
$ws.Cells.Item(31, 5).Value = 6750
$ws.Cells.Item(31, 7).Value = '10-Nov-2025 00:00:00'
$ws.Cells.Item(31, 22).Value = 'RS'
$ws.Cells.Item(31, 27).Value = 5
$ws.Cells.Item(31, 28).Value = 0.9
$ws.Cells.Item(31, 31).Value = 5.9

# --- Add new rows 32-38 ---
# Row 32
$ws.Cells.Item(32, 1).Value = 16104
$ws.Cells.Item(32, 2).Value = 'KARRI ABHINAY CHARVIK'
$ws.Cells.Item(32, 3).Value = 9177965994
$ws.Cells.Item(32, 4).Value = 8350
$ws.Cells.Item(32, 5).Value = 8350
$ws.Cells.Item(32, 6).Value = '10-Nov-2025 09:27:15'
$ws.Cells.Item(32, 7).Value = '11-Nov-2025 00:00:00'
$ws.Cells.Item(32, 8).Value = 'TRANSACTION IS SUCCESSFUL'
$ws.Cells.Item(32, 9).Value = 'OK'
$ws.Cells.Item(32, 10).Value = 'lVl'
$ws.Cells.Item(32, 11).Value = 'eight thousand three hundred fifty'
$ws.Cells.Item(32, 12).Value = 100000036600
$ws.Cells.Item(32, 13).Value = 'SALESIAN EDUCATION SOCIETY'
$ws.Cells.Item(32, 14).Value = 753702
$ws.Cells.Item(32, 15).Value = 1234
$ws.Cells.Item(32, 16).Value = 11000315754150
$ws.Cells.Item(32, 17).Value = 1762746971
$ws.Cells.Item(32, 18).Value = 531445403467
$ws.Cells.Item(32, 19).Value = 'INR'
$ws.Cells.Item(32, 20).Value = 'sale'
$ws.Cells.Item(32, 21).Value = 'ICICI UPI QR'
$ws.Cells.Item(32, 22).Value = 'RS'
$ws.Cells.Item(32, 23).Value = 'SIBL0000899'
$ws.Cells.Item(32, 24).Value = 'MERCHANT'
$ws.Cells.Item(32, 25).Value = 'UPI'
$ws.Cells.Item(32, 26).Value = 'kotakschoolvsp@gmail.com'
$ws.Cells.Item(32, 27).Value = 5
$ws.Cells.Item(32, 28).Value = 0.9
$ws.Cells.Item(32, 29).Value = 0
$ws.Cells.Item(32, 30).Value = 0
$ws.Cells.Item(32, 31).Value = 5.9
$ws.Cells.Item(32, 32).Value = 'KOTAK SALESIAN PRIMARY SCHOOL'
$ws.Cells.Item(32, 34).Value = 'REGULAR'
$ws.Cells.Item(32, 35).Value = 18807
$c = $ws.Cells.Item(32, 36)
$c.NumberFormat = "@"
$c.Value = '265833'
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 37)
$c.NumberFormat = "@"
$c.Value = '2050'
$c.Style = "Normal"

# Row 33
$ws.Cells.Item(33, 1).Value = 16052
$ws.Cells.Item(33, 2).Value = 'GURRALA RUSHABH NARAYAN'
$ws.Cells.Item(33, 3).Value = 9704995001
$ws.Cells.Item(33, 4).Value = 10750
$ws.Cells.Item(33, 5).Value = 10750
$ws.Cells.Item(33, 6).Value = '10-Nov-2025 21:17:01'
$ws.Cells.Item(33, 8).Value = 'TRANSACTION IS SUCCESSFUL'
$ws.Cells.Item(33, 9).Value = 'OK'
$ws.Cells.Item(33, 10).Value = 'VlllX'
$ws.Cells.Item(33, 11).Value = 'ten thousand seven hundred fifty'
$ws.Cells.Item(33, 12).Value = 100000036600
$ws.Cells.Item(33, 13).Value = 'SALESIAN EDUCATION SOCIETY'
$ws.Cells.Item(33, 14).Value = 753702
$ws.Cells.Item(33, 15).Value = 1234
$ws.Cells.Item(33, 16).Value = 11000316033651
$ws.Cells.Item(33, 17).Value = 1762789270
$ws.Cells.Item(33, 18).Value = 531449495645
$ws.Cells.Item(33, 19).Value = 'INR'
$ws.Cells.Item(33, 20).Value = 'sale'
$ws.Cells.Item(33, 21).Value = 'ICICI UPI QR'
$ws.Cells.Item(33, 22).Value = 'RNS'
$ws.Cells.Item(33, 23).Value = 'SIBL0000899'
$ws.Cells.Item(33, 24).Value = 'MERCHANT'
$ws.Cells.Item(33, 25).Value = 'UPI'
$ws.Cells.Item(33, 26).Value = 'kotakschoolvsp@gmail.com'
$ws.Cells.Item(33, 27).Value = 5
$ws.Cells.Item(33, 28).Value = 0.9
$ws.Cells.Item(33, 29).Value = 0
$ws.Cells.Item(33, 30).Value = 0
$ws.Cells.Item(33, 31).Value = 5.9
$ws.Cells.Item(33, 32).Value = 'KOTAK SALESIAN SECONDARY SCHOOL ICSE'
$ws.Cells.Item(33, 34).Value = 'REGULAR'
$ws.Cells.Item(33, 35).Value = 19325
$c = $ws.Cells.Item(33, 36)
$c.NumberFormat = "@"
$c.Value = '266322'
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 37)
$c.NumberFormat = "@"
$c.Value = '2052'
$c.Style = "Normal"

# Row 34
$ws.Cells.Item(34, 1).Value = 17315
$ws.Cells.Item(34, 2).Value = 'VADIGI DIVYESH SAI PAVAN'
$ws.Cells.Item(34, 3).Value = 8367352524
$ws.Cells.Item(34, 4).Value = 6750
$ws.Cells.Item(34, 6).Value = '11-Nov-2025 05:53:50'
$ws.Cells.Item(34, 8).Value = 'TRANSACTION IS SUCCESSFUL'
$ws.Cells.Item(34, 9).Value = 'OK'
$ws.Cells.Item(34, 10).Value = 'PREKGUKG'
$ws.Cells.Item(34, 11).Value = 'six thousand seven hundred fifty'
$ws.Cells.Item(34, 12).Value = 100000036600
$ws.Cells.Item(34, 13).Value = 'SALESIAN EDUCATION SOCIETY'
$ws.Cells.Item(34, 14).Value = 753702
$ws.Cells.Item(34, 15).Value = 1234
$ws.Cells.Item(34, 16).Value = 11000316067938
$ws.Cells.Item(34, 17).Value = 1762821403
$ws.Cells.Item(34, 18).Value = 30020732085
$ws.Cells.Item(34, 19).Value = 'INR'
$ws.Cells.Item(34, 20).Value = 'sale'
$ws.Cells.Item(34, 21).Value = 'ICICI UPI QR'
$ws.Cells.Item(34, 22).Value = 'NRNS'
$ws.Cells.Item(34, 23).Value = 'SIBL0000899'
$ws.Cells.Item(34, 24).Value = 'MERCHANT'
$ws.Cells.Item(34, 25).Value = 'UPI'
$ws.Cells.Item(34, 26).Value = 'kotakschoolvsp@gmail.com'
$ws.Cells.Item(34, 29).Value = 0
$ws.Cells.Item(34, 30).Value = 0
$ws.Cells.Item(34, 32).Value = 'KOTAK SALESIAN SCHOOL MANAGEMENT ACCOUNT'
$ws.Cells.Item(34, 34).Value = 'REGULAR'
$ws.Cells.Item(34, 35).Value = 20023
$c = $ws.Cells.Item(34, 36)
$c.NumberFormat = "@"
$c.Value = '265150'
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 37)
$c.NumberFormat = "@"
$c.Value = '2047'
$c.Style = "Normal"
$ws.Cells.Item(34, 38).Value = 'UPI INTENT'

# Row 35
$ws.Cells.Item(35, 1).Value = 16733
$ws.Cells.Item(35, 2).Value = 'HARINI CHAND REDDY PILAKA'
$ws.Cells.Item(35, 3).Value = 9866815664
$ws.Cells.Item(35, 4).Value = 8350
$ws.Cells.Item(35, 6).Value = '11-Nov-2025 10:09:11'
$ws.Cells.Item(35, 8).Value = 'TRANSACTION IS SUCCESSFUL'
$ws.Cells.Item(35, 9).Value = 'OK'
$ws.Cells.Item(35, 10).Value = 'lVl'
$ws.Cells.Item(35, 11).Value = 'eight thousand three hundred fifty'
$ws.Cells.Item(35, 12).Value = 100000036600
$ws.Cells.Item(35, 13).Value = 'SALESIAN EDUCATION SOCIETY'
$ws.Cells.Item(35, 14).Value = 753702
$ws.Cells.Item(35, 15).Value = 1234
$ws.Cells.Item(35, 16).Value = 11000316099462
$ws.Cells.Item(35, 17).Value = 1762835725
$ws.Cells.Item(35, 18).Value = 108551696284
$ws.Cells.Item(35, 19).Value = 'INR'
$ws.Cells.Item(35, 20).Value = 'sale'
$ws.Cells.Item(35, 21).Value = 'ICICI UPI QR'
$ws.Cells.Item(35, 22).Value = 'NRNS'
$ws.Cells.Item(35, 23).Value = 'SIBL0000899'
$ws.Cells.Item(35, 24).Value = 'MERCHANT'
$ws.Cells.Item(35, 25).Value = 'UPI'
$ws.Cells.Item(35, 26).Value = 'kotakschoolvsp@gmail.com'
$ws.Cells.Item(35, 29).Value = 0
$ws.Cells.Item(35, 30).Value = 0
$ws.Cells.Item(35, 32).Value = 'KOTAK SALESIAN PRIMARY SCHOOL'
$ws.Cells.Item(35, 34).Value = 'REGULAR'
$ws.Cells.Item(35, 35).Value = 18753
$c = $ws.Cells.Item(35, 36)
$c.NumberFormat = "@"
$c.Value = '265774'
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 37)
$c.NumberFormat = "@"
$c.Value = '2050'
$c.Style = "Normal"

# Row 36
$ws.Cells.Item(36, 1).Value = 16734
$ws.Cells.Item(36, 2).Value = 'HASINI CHAND REDDY PILAKA'
$ws.Cells.Item(36, 3).Value = 9866815664
$ws.Cells.Item(36, 4).Value = 11350
$ws.Cells.Item(36, 6).Value = '11-Nov-2025 10:10:40'
$ws.Cells.Item(36, 8).Value = 'TRANSACTION IS SUCCESSFUL'
$ws.Cells.Item(36, 9).Value = 'OK'
$ws.Cells.Item(36, 10).Value = 'VlllX'
$ws.Cells.Item(36, 11).Value = 'eleven thousand three hundred fifty'
$ws.Cells.Item(36, 12).Value = 100000036600
$ws.Cells.Item(36, 13).Value = 'SALESIAN EDUCATION SOCIETY'
$ws.Cells.Item(36, 14).Value = 753702
$ws.Cells.Item(36, 15).Value = 1234
$ws.Cells.Item(36, 16).Value = 11000316100223
$ws.Cells.Item(36, 17).Value = 1762836014
$ws.Cells.Item(36, 18).Value = 108551703754
$ws.Cells.Item(36, 19).Value = 'INR'
$ws.Cells.Item(36, 20).Value = 'sale'
$ws.Cells.Item(36, 21).Value = 'ICICI UPI QR'
$ws.Cells.Item(36, 22).Value = 'NRNS'
$ws.Cells.Item(36, 23).Value = 'SIBL0000899'
$ws.Cells.Item(36, 24).Value = 'MERCHANT'
$ws.Cells.Item(36, 25).Value = 'UPI'
$ws.Cells.Item(36, 26).Value = 'kotakschoolvsp@gmail.com'
$ws.Cells.Item(36, 29).Value = 0
$ws.Cells.Item(36, 30).Value = 0
$ws.Cells.Item(36, 32).Value = 'KOTAK SALESIAN SECONDARY SCHOOL ICSE'
$ws.Cells.Item(36, 34).Value = 'REGULAR'
$ws.Cells.Item(36, 35).Value = 19425
$c = $ws.Cells.Item(36, 36)
$c.NumberFormat = "@"
$c.Value = '266398'
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 37)
$c.NumberFormat = "@"
$c.Value = '2053'
$c.Style = "Normal"

# Row 37
$ws.Cells.Item(37, 1).Value = 16927
$ws.Cells.Item(37, 2).Value = 'THUTTA DHANASHVI'
$ws.Cells.Item(37, 3).Value = 7337239208
$ws.Cells.Item(37, 4).Value = 7750
$ws.Cells.Item(37, 6).Value = '11-Nov-2025 16:21:56'
$ws.Cells.Item(37, 8).Value = 'TRANSACTION IS SUCCESSFUL'
$ws.Cells.Item(37, 9).Value = 'OK'
$ws.Cells.Item(37, 10).Value = 'lVl'
$ws.Cells.Item(37, 11).Value = 'seven thousand seven hundred fifty'
$ws.Cells.Item(37, 12).Value = 100000036600
$ws.Cells.Item(37, 13).Value = 'SALESIAN EDUCATION SOCIETY'
$ws.Cells.Item(37, 14).Value = 753702
$ws.Cells.Item(37, 15).Value = 1234
$ws.Cells.Item(37, 16).Value = 11000316191666
$ws.Cells.Item(37, 17).Value = 1762858305
$ws.Cells.Item(37, 18).Value = 783823510579
$ws.Cells.Item(37, 19).Value = 'INR'
$ws.Cells.Item(37, 20).Value = 'sale'
$ws.Cells.Item(37, 21).Value = 'ICICI UPI QR'
$ws.Cells.Item(37, 22).Value = 'NRNS'
$ws.Cells.Item(37, 23).Value = 'SIBL0000899'
$ws.Cells.Item(37, 24).Value = 'MERCHANT'
$ws.Cells.Item(37, 25).Value = 'UPI'
$ws.Cells.Item(37, 26).Value = 'kotakschoolvsp@gmail.com'
$ws.Cells.Item(37, 29).Value = 0
$ws.Cells.Item(37, 30).Value = 0
$ws.Cells.Item(37, 32).Value = 'KOTAK SALESIAN PRIMARY SCHOOL'
$ws.Cells.Item(37, 34).Value = 'REGULAR'
$ws.Cells.Item(37, 35).Value = 18339
$c = $ws.Cells.Item(37, 36)
$c.NumberFormat = "@"
$c.Value = '265347'
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 37)
$c.NumberFormat = "@"
$c.Value = '2048'
$c.Style = "Normal"
$ws.Cells.Item(37, 38).Value = 'UPI INTENT'

# Row 38
$ws.Cells.Item(38, 1).Value = 17177
$ws.Cells.Item(38, 2).Value = 'THUTTA HARI CHARANI'
$ws.Cells.Item(38, 3).Value = 7337239208
$ws.Cells.Item(38, 4).Value = 6750
$ws.Cells.Item(38, 6).Value = '11-Nov-2025 16:21:09'
$ws.Cells.Item(38, 8).Value = 'TRANSACTION IS SUCCESSFUL'
$ws.Cells.Item(38, 9).Value = 'OK'
$ws.Cells.Item(38, 10).Value = 'PREKGUKG'
$ws.Cells.Item(38, 11).Value = 'six thousand seven hundred fifty'
$ws.Cells.Item(38, 12).Value = 100000036600
$ws.Cells.Item(38, 13).Value = 'SALESIAN EDUCATION SOCIETY'
$ws.Cells.Item(38, 14).Value = 753702
$ws.Cells.Item(38, 15).Value = 1234
$ws.Cells.Item(38, 16).Value = 11000316192344
$ws.Cells.Item(38, 17).Value = 1762858254
$ws.Cells.Item(38, 18).Value = 704447957133
$ws.Cells.Item(38, 19).Value = 'INR'
$ws.Cells.Item(38, 20).Value = 'sale'
$ws.Cells.Item(38, 21).Value = 'ICICI UPI QR'
$ws.Cells.Item(38, 22).Value = 'NRNS'
$ws.Cells.Item(38, 23).Value = 'SIBL0000899'
$ws.Cells.Item(38, 24).Value = 'MERCHANT'
$ws.Cells.Item(38, 25).Value = 'UPI'
$ws.Cells.Item(38, 26).Value = 'kotakschoolvsp@gmail.com'
$ws.Cells.Item(38, 29).Value = 0
$ws.Cells.Item(38, 30).Value = 0
$ws.Cells.Item(38, 32).Value = 'KOTAK SALESIAN SCHOOL MANAGEMENT ACCOUNT'
$ws.Cells.Item(38, 34).Value = 'REGULAR'
$ws.Cells.Item(38, 35).Value = 19892
$c = $ws.Cells.Item(38, 36)
$c.NumberFormat = "@"
$c.Value = '265011'
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 37)
$c.NumberFormat = "@"
$c.Value = '2047'
$c.Style = "Normal"
$ws.Cells.Item(38, 38).Value = 'UPI INTENT'

Write-Host "Edit complete"